# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 341
    "F4"  = 10536
    "F6"  = 959
    "F7"  = 71
    "F8"  = 1295
    "F9"  = 7793
    "F15" = 3231
    "F17" = 320
    "F18" = 714
    "F20" = 1048
    "F21" = 281
    "F22" = 90
    "F23" = 1679
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
